# templateEstudiantesMaes.xlsx - turn the fixed column-header template into a
# blank (content-less) layout, and tweak the print/view setup, matching the
# "preparado por fin un funcional de creacion de excel y pdf" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The B10:B23 column used to hold the fixed field labels (CODIGO, NOMBRE,
# CORREO, CEDULA, CIUDAD, CONDICION, NIVEL, DIRECTOR, CODIRECTOR 1,
# CODIRECTOR 2, REGLAMENTO, TEMA, FECHA, CONCEPTO). Clear the text but keep
# the cell formatting/styles untouched so the table skeleton stays usable.
$ws.Range("B10:B23").ClearContents()

# Re-orient the worksheet for printing (landscape instead of portrait) and
# adjust the page margins (left 0.6cm, right 14.6cm, top/bottom 1.9cm,
# header/footer 0.8cm -- expressed in points as Excel's object model wants).
$ps = $ws.PageSetup
$ps.Orientation = 2
$ps.LeftMargin = 17.007874015748033
$ps.RightMargin = 413.85826771653547
$ps.TopMargin = 53.85826771653544
$ps.BottomMargin = 53.85826771653544
$ps.HeaderMargin = 22.677165354330707
$ps.FooterMargin = 22.677165354330707

# Move the active selection down to where the user was last working.
$ws.Activate()
$ws.Range("D15").Select()
